$d = $word.ActiveDocument

# The document contains five "<id>p093r_aN</id>" markers (N = 1..5), each
# split across three runs: "<id>", "p093r_aN", "</id>". The edit collapses
# each of these into a single run reading "<id>p093r_N</id>" (dropping the
# "a" from the id and merging the three runs/formats into one, keeping the
# Courier-New formatting that the "<id>" / "</id>" runs already had).

for ($i = 1; $i -le 5; $i++) {
    $old = "<id>p093r_a$i</id>"
    $new = "<id>p093r_$i</id>"
    $found = $d.Content.Find.Execute($old, $false, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    Write-Host "Replaced #$i : $found"
}
